# Auto-generated Excel COM-interop script
# Applies updated H..N financial figures to multiple sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 2500  # H32
$ws.Cells.Item(32, 9).Value = 2500  # I32
$ws.Cells.Item(32, 11).Value = 2500  # K32
$ws.Cells.Item(32, 13).Value = -2174  # M32
$ws.Cells.Item(40, 8).Value = 13158.333  # H40
$ws.Cells.Item(40, 9).Value = 11737.5  # I40
$ws.Cells.Item(40, 11).Value = 11737.5  # K40
$ws.Cells.Item(40, 13).Value = -11562.5  # M40
$ws.Cells.Item(55, 8).Value = 509.94446  # H55
$ws.Cells.Item(55, 10).Value = 844.7778  # J55
$ws.Cells.Item(55, 12).Value = 844.7778  # L55
$ws.Cells.Item(55, 14).Value = -1272.7778  # N55
$ws.Cells.Item(125, 8).Value = 958.6667  # H125
$ws.Cells.Item(125, 9).Value = 841.4  # I125
$ws.Cells.Item(125, 10).Value = 1105.25  # J125
$ws.Cells.Item(125, 11).Value = 7572.599999999999  # K125
$ws.Cells.Item(125, 12).Value = 9947.25  # L125
$ws.Cells.Item(125, 13).Value = -5112.599999999999  # M125
$ws.Cells.Item(125, 14).Value = -14867.25  # N125
$ws.Cells.Item(135, 8).Value = 1680.8  # H135
$ws.Cells.Item(135, 9).Value = 1680.8  # I135
$ws.Cells.Item(135, 11).Value = 15127.2  # K135
$ws.Cells.Item(135, 13).Value = -12592.2  # M135
$ws.Cells.Item(138, 8).Value = 2933.8308  # H138
$ws.Cells.Item(138, 10).Value = 2931.6724  # J138
$ws.Cells.Item(138, 12).Value = 8795.0172  # L138
$ws.Cells.Item(138, 14).Value = -19075.0172  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(16, 8).Value = 1849  # H16
$ws.Cells.Item(16, 9).Value = 1918.8  # I16
$ws.Cells.Item(16, 10).Value = 1500  # J16
$ws.Cells.Item(16, 11).Value = 1918.8  # K16
$ws.Cells.Item(16, 12).Value = 1500  # L16
$ws.Cells.Item(16, 13).Value = -1631.8  # M16
$ws.Cells.Item(16, 14).Value = -2074  # N16
$ws.Cells.Item(32, 8).Value = 12292.16  # H32
$ws.Cells.Item(32, 9).Value = 5900.5933  # I32
$ws.Cells.Item(32, 11).Value = 5900.5933  # K32
$ws.Cells.Item(32, 13).Value = -5613.5933  # M32
$ws.Cells.Item(45, 8).Value = 2870.5  # H45
$ws.Cells.Item(45, 9).Value = 2811  # I45
$ws.Cells.Item(45, 11).Value = 2811  # K45
$ws.Cells.Item(45, 13).Value = -2434  # M45
$ws.Cells.Item(74, 8).Value = 1615.5555  # H74
$ws.Cells.Item(74, 9).Value = 1615.5555  # I74
$ws.Cells.Item(74, 11).Value = 1615.5555  # K74
$ws.Cells.Item(74, 13).Value = -741.5554999999999  # M74
$ws.Cells.Item(77, 8).Value = 1615.5555  # H77
$ws.Cells.Item(77, 9).Value = 1615.5555  # I77
$ws.Cells.Item(77, 11).Value = 8077.7775  # K77
$ws.Cells.Item(77, 13).Value = -3709.7775  # M77
$ws.Cells.Item(110, 8).Value = 2602.6428  # H110
$ws.Cells.Item(110, 9).Value = 2149.0908  # I110
$ws.Cells.Item(110, 11).Value = 2149.0908  # K110
$ws.Cells.Item(110, 13).Value = -104.0907999999999  # M110
$ws.Cells.Item(122, 8).Value = 3233.1516  # H122
$ws.Cells.Item(122, 9).Value = 2414.889  # I122
$ws.Cells.Item(122, 11).Value = 7244.667  # K122
$ws.Cells.Item(122, 13).Value = -4794.667  # M122
$ws.Cells.Item(132, 8).Value = 2111.5715  # H132
$ws.Cells.Item(132, 9).Value = 1963.6666  # I132
$ws.Cells.Item(132, 11).Value = 5890.9998  # K132
$ws.Cells.Item(132, 13).Value = -3360.9998  # M132

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(15, 8).Value = 20000  # H15
$ws.Cells.Item(15, 10).Value = 20000  # J15
$ws.Cells.Item(15, 12).Value = 20000  # L15
$ws.Cells.Item(15, 14).Value = -20454  # N15
$ws.Cells.Item(16, 8).Value = 12498.286  # H16
$ws.Cells.Item(16, 9).Value = 12498.286  # I16
$ws.Cells.Item(16, 11).Value = 12498.286  # K16
$ws.Cells.Item(16, 13).Value = -12328.286  # M16
$ws.Cells.Item(22, 8).Value = 748.75  # H22
$ws.Cells.Item(22, 9).Value = 331.83334  # I22
$ws.Cells.Item(22, 11).Value = 331.83334  # K22
$ws.Cells.Item(22, 13).Value = -158.83334  # M22
$ws.Cells.Item(33, 8).Value = 1021  # H33
$ws.Cells.Item(33, 9).Value = 1021  # I33
$ws.Cells.Item(33, 11).Value = 1021  # K33
$ws.Cells.Item(33, 13).Value = -685  # M33
$ws.Cells.Item(52, 8).Value = 59999.5  # H52
$ws.Cells.Item(52, 10).Value = 59999.5  # J52
$ws.Cells.Item(52, 12).Value = 59999.5  # L52
$ws.Cells.Item(52, 14).Value = -60525.5  # N52
$ws.Cells.Item(107, 8).Value = 2392.2856  # H107
$ws.Cells.Item(107, 9).Value = 2391  # I107
$ws.Cells.Item(107, 11).Value = 2391  # K107
$ws.Cells.Item(107, 13).Value = -471  # M107
$ws.Cells.Item(121, 8).Value = 59999.5  # H121
$ws.Cells.Item(121, 10).Value = 59999.5  # J121
$ws.Cells.Item(121, 12).Value = 59999.5  # L121
$ws.Cells.Item(121, 14).Value = -63493.5  # N121
$ws.Cells.Item(131, 8).Value = 69999  # H131
$ws.Cells.Item(131, 10).Value = 69999  # J131
$ws.Cells.Item(131, 12).Value = 69999  # L131
$ws.Cells.Item(131, 14).Value = -80079  # N131

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6652.479  # H31
$ws.Cells.Item(31, 10).Value = 13734.267  # J31
$ws.Cells.Item(31, 12).Value = 13734.267  # L31
$ws.Cells.Item(31, 14).Value = -14324.267  # N31
$ws.Cells.Item(34, 8).Value = 6652.479  # H34
$ws.Cells.Item(34, 10).Value = 13734.267  # J34
$ws.Cells.Item(34, 12).Value = 13734.267  # L34
$ws.Cells.Item(34, 14).Value = -14138.267  # N34
$ws.Cells.Item(41, 8).Value = 29999  # H41
$ws.Cells.Item(41, 9).Value = 0  # I41
$ws.Cells.Item(41, 10).Value = 29999  # J41
$ws.Cells.Item(41, 11).Value = 0  # K41
$ws.Cells.Item(41, 12).Value = 29999  # L41
$ws.Cells.Item(41, 13).ClearContents()  # M41
$ws.Cells.Item(41, 14).Value = -30855  # N41
$ws.Cells.Item(132, 8).Value = 2426.3215  # H132
$ws.Cells.Item(132, 9).Value = 2210.9092  # I132
$ws.Cells.Item(132, 11).Value = 6632.7276  # K132
$ws.Cells.Item(132, 13).Value = -4102.7276  # M132

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(114, 8).Value = 1387  # H114
$ws.Cells.Item(114, 9).Value = 1617.3  # I114
$ws.Cells.Item(114, 10).Value = 1003.1667  # J114
$ws.Cells.Item(114, 11).Value = 4851.9  # K114
$ws.Cells.Item(114, 12).Value = 3009.5001  # L114
$ws.Cells.Item(114, 13).Value = -1597.9  # M114
$ws.Cells.Item(114, 14).Value = -9517.500100000001  # N114

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 7937.25  # H70
$ws.Cells.Item(70, 9).Value = 6999.5  # I70
$ws.Cells.Item(70, 11).Value = 6999.5  # K70
$ws.Cells.Item(70, 13).Value = -6729.5  # M70
$ws.Cells.Item(73, 8).Value = 7937.25  # H73
$ws.Cells.Item(73, 9).Value = 6999.5  # I73
$ws.Cells.Item(73, 11).Value = 6999.5  # K73
$ws.Cells.Item(73, 13).Value = -6063.5  # M73
$ws.Cells.Item(97, 8).Value = 889.9  # H97
$ws.Cells.Item(97, 9).Value = 1047.0667  # I97
$ws.Cells.Item(97, 10).Value = 418.4  # J97
$ws.Cells.Item(97, 11).Value = 1047.0667  # K97
$ws.Cells.Item(97, 12).Value = 418.4  # L97
$ws.Cells.Item(97, 13).Value = -551.0667000000001  # M97
$ws.Cells.Item(97, 14).Value = -1410.4  # N97
$ws.Cells.Item(99, 8).Value = 7226.25  # H99
$ws.Cells.Item(99, 9).Value = 7226.25  # I99
$ws.Cells.Item(99, 11).Value = 7226.25  # K99
$ws.Cells.Item(99, 13).Value = -4980.25  # M99
$ws.Cells.Item(102, 8).Value = 4214.857  # H102
$ws.Cells.Item(102, 9).Value = 3480.7334  # I102
$ws.Cells.Item(102, 10).Value = 6050.1665  # J102
$ws.Cells.Item(102, 11).Value = 3480.7334  # K102
$ws.Cells.Item(102, 12).Value = 6050.1665  # L102
$ws.Cells.Item(102, 13).Value = -1858.7334  # M102
$ws.Cells.Item(102, 14).Value = -9294.1665  # N102
$ws.Cells.Item(117, 8).Value = 66310  # H117
$ws.Cells.Item(117, 10).Value = 66310  # J117
$ws.Cells.Item(117, 12).Value = 66310  # L117
$ws.Cells.Item(117, 14).Value = -73194  # N117

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(29, 8).Value = 44000  # H29
$ws.Cells.Item(29, 10).Value = 44000  # J29
$ws.Cells.Item(29, 12).Value = 44000  # L29
$ws.Cells.Item(29, 14).Value = -44590  # N29
$ws.Cells.Item(40, 8).Value = 5752.9414  # H40
$ws.Cells.Item(40, 9).Value = 4753.846  # I40
$ws.Cells.Item(40, 10).Value = 9000  # J40
$ws.Cells.Item(40, 11).Value = 4753.846  # K40
$ws.Cells.Item(40, 12).Value = 9000  # L40
$ws.Cells.Item(40, 13).Value = -4617.846  # M40
$ws.Cells.Item(40, 14).Value = -9272  # N40
$ws.Cells.Item(116, 8).Value = 140000  # H116
$ws.Cells.Item(116, 10).Value = 140000  # J116
$ws.Cells.Item(116, 12).Value = 140000  # L116
$ws.Cells.Item(116, 14).Value = -149178  # N116
$ws.Cells.Item(122, 8).Value = 3394.1428  # H122
$ws.Cells.Item(122, 9).Value = 2675.25  # I122
$ws.Cells.Item(122, 10).Value = 3486.9033  # J122
$ws.Cells.Item(122, 11).Value = 8025.75  # K122
$ws.Cells.Item(122, 12).Value = 10460.7099  # L122
$ws.Cells.Item(122, 13).Value = -5575.75  # M122
$ws.Cells.Item(122, 14).Value = -15360.7099  # N122

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 4080.4167  # H81
$ws.Cells.Item(81, 9).Value = 2110.889  # I81
$ws.Cells.Item(81, 11).Value = 4221.778  # K81
$ws.Cells.Item(81, 13).Value = -3160.778  # M81
$ws.Cells.Item(84, 8).Value = 4080.4167  # H84
$ws.Cells.Item(84, 9).Value = 2110.889  # I84
$ws.Cells.Item(84, 11).Value = 21108.89  # K84
$ws.Cells.Item(84, 13).Value = -15804.89  # M84
$ws.Cells.Item(107, 8).Value = 718.9375  # H107
$ws.Cells.Item(107, 9).Value = 379.41666  # I107
$ws.Cells.Item(107, 11).Value = 1138.24998  # K107
$ws.Cells.Item(107, 13).Value = 781.75002  # M107
$ws.Cells.Item(132, 8).Value = 1991.9762  # H132
$ws.Cells.Item(132, 9).Value = 1820.9459  # I132
$ws.Cells.Item(132, 10).Value = 3257.6  # J132
$ws.Cells.Item(132, 11).Value = 5462.8377  # K132
$ws.Cells.Item(132, 12).Value = 9772.8  # L132
$ws.Cells.Item(132, 13).Value = -2932.8377  # M132
$ws.Cells.Item(132, 14).Value = -14832.8  # N132
$ws.Cells.Item(136, 8).Value = 3596.111  # H136
$ws.Cells.Item(136, 9).Value = 3852.1428  # I136
$ws.Cells.Item(136, 11).Value = 11556.4284  # K136
$ws.Cells.Item(136, 13).Value = -9006.4284  # M136

